# Applies the cryptos list refresh described in the commit:
# "Updated cryptos list on Mon Jan 22 15:56:04 UTC 2024 with GitHub Actions"
#
# All data cells on the sheet are stored as text (coin name, link, price,
# 1h-volume-change %). Prices such as "2.730.97" or "0.110" must stay text
# -- Excel's COM Value setter auto-coerces numeric-looking strings to
# numbers, which would both change the cell type and silently drop
# formatting-significant trailing zeros. Prefixing with an apostrophe
# (the same 'force text' marker Excel uses for manual entry) keeps every
# written value a literal string, exactly like the source cells.
function Set-Text($row, $col, $text) {
    $ws.Cells.Item($row, $col).Value = "'" + $text
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
Set-Text 2 4 '40.678.20'
Set-Text 2 5 '  -2.31%  '

# Row 3
Set-Text 3 4 '2.374.86'
Set-Text 3 5 '  -3.77%  '

# Row 4
Set-Text 4 5 '  -0.09%  '

# Row 5
Set-Text 5 4 '310.64'
Set-Text 5 5 '  -2.39%  '

# Row 6
Set-Text 6 4 '87.04'
Set-Text 6 5 '  -5.91%  '

# Row 7
Set-Text 7 4 '0.528'
Set-Text 7 5 '  -4.20%  '

# Row 8
Set-Text 8 5 '  -0.01%  '

# Row 9
Set-Text 9 4 '0.494'
Set-Text 9 5 '  -4.02%  '

# Row 10
Set-Text 10 4 '0.0842'
Set-Text 10 5 '  -3.33%  '

# Row 11
Set-Text 11 4 '30.58'
Set-Text 11 5 '  -7.32%  '

# Row 12
Set-Text 12 4 '0.110'
Set-Text 12 5 '  -0.76%  '

# Row 13
Set-Text 13 4 '2.730.97'
Set-Text 13 5 '  -4.10%  '

# Row 14
Set-Text 14 4 '6.54'
Set-Text 14 5 '  -4.98%  '

# Row 15
Set-Text 15 4 '15.05'
Set-Text 15 5 '  -2.95%  '

# Row 16
Set-Text 16 4 '2.373.89'
Set-Text 16 5 '  -4.46%  '

# Row 17
Set-Text 17 4 '0.760'
Set-Text 17 5 '  -4.51%  '

# Row 18
Set-Text 18 4 '40.494.18'
Set-Text 18 5 '  -2.63%  '

# Row 19
Set-Text 19 4 '0.0₃0911'
Set-Text 19 5 '  -3.77%  '

# Row 20
Set-Text 20 4 '6.13'
Set-Text 20 5 '  -4.90%  '

# Row 21
Set-Text 21 4 '68.54'
Set-Text 21 5 '  -3.28%  '

# Row 22
Set-Text 22 4 '10.74'
Set-Text 22 5 '  -4.95%  '

# Row 23
Set-Text 23 4 '235.13'
Set-Text 23 5 '  -2.36%  '

# Row 24
Set-Text 24 4 '2.59'
Set-Text 24 5 '  -5.80%  '

# Row 25
Set-Text 25 5 '  +0.14%  '

# Row 26
Set-Text 26 5 '  -8.16%  '

# Row 27
Set-Text 27 4 '23.67'
Set-Text 27 5 '  -4.95%  '

# Row 28
Set-Text 28 4 '2.20'
Set-Text 28 5 '  -2.12%  '

# Row 29
Set-Text 29 4 '9.30'
Set-Text 29 5 '  -4.22%  '

# Row 30
Set-Text 30 4 '33.82'
Set-Text 30 5 '  -7.61%  '

# Row 31
Set-Text 31 4 '152.59'
Set-Text 31 5 '  -3.42%  '

# Row 32
Set-Text 32 5 '  -0.08%  '

# Row 33
Set-Text 33 4 '5.22'
Set-Text 33 5 '  -4.81%  '

# Row 34
Set-Text 34 4 '0.0729'
Set-Text 34 5 '  -4.22%  '

# Row 35
Set-Text 35 5 '  -5.33%  '

# Row 36
Set-Text 36 5 '  -1.95%  '

# Row 37
Set-Text 37 2 'Celestia'
Set-Text 37 3 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
Set-Text 37 4 '15.87'
Set-Text 37 5 '  -8.31%  '

# Row 38
Set-Text 38 2 'Kaspa'
Set-Text 38 3 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-Text 38 4 '0.0996'
Set-Text 38 5 '  -4.44%  '

# Row 39
Set-Text 39 4 '2.75'
Set-Text 39 5 '  -5.49%  '

# Row 40
Set-Text 40 4 '1.71'
Set-Text 40 5 '  -7.67%  '

# Row 41
Set-Text 41 4 '3.83'
Set-Text 41 5 '  -4.17%  '

# Row 42
Set-Text 42 4 '2.38'
Set-Text 42 5 '  -6.22%  '

# Row 43
Set-Text 43 4 '1.961.48'
Set-Text 43 5 '  -1.42%  '

# Row 44
Set-Text 44 4 '0.0269'
Set-Text 44 5 '  -5.05%  '

# Row 45
Set-Text 45 5 '  -7.37%  '

# Row 46
Set-Text 46 4 '9.45'
Set-Text 46 5 '  +0.44%  '

# Row 47
Set-Text 47 4 '2.70'
Set-Text 47 5 '  -8.88%  '

# Row 48
Set-Text 48 4 '2.594.08'
Set-Text 48 5 '  -4.15%  '

# Row 49
Set-Text 49 4 '93.30'
Set-Text 49 5 '  -4.40%  '

# Row 50
Set-Text 50 4 '72.36'
Set-Text 50 5 '  -4.42%  '

# Row 51
Set-Text 51 4 '50.46'
Set-Text 51 5 '  -4.05%  '

